$rows = @(
    @{ Row=2; A=90; C='DORO CAT 豆乳貓 活性碳與小蘇打雙重除臭顆粒低過敏極細豆腐貓砂, 無香, 7L, 2包'; D='$298'; E='($2.13/100ml)'; F='DORO CAT 豆乳貓 活性碳與小蘇打雙重除臭顆粒低過敏極細豆腐貓砂, 無香, 7L, 2包 特價 6折 $500 $298 ($2.13/100ml) 7折 優惠券 缺貨 免運 滿 $490 ( 52 )'; G=45789.03155092592 }
    @{ Row=3; A=89; C='DORO CAT 豆乳貓 豆腐貓砂, 原味, 7L, 6袋'; D='$675'; E='($1.61/100ml)'; F='DORO CAT 豆乳貓 豆腐貓砂, 原味, 7L, 6袋 特價 45折 $1,500 $675 ($1.61/100ml) 7折 優惠券 缺貨 免運 ( 489 )'; G=45789.03155092592 }
    @{ Row=4; A=88; C='DORO CAT 豆乳貓 活性碳與小蘇打雙重除臭顆粒低過敏極細豆腐貓砂, 無香, 7L, 6包'; D='$675'; E='($1.61/100ml)'; F='DORO CAT 豆乳貓 活性碳與小蘇打雙重除臭顆粒低過敏極細豆腐貓砂, 無香, 7L, 6包 特價 45折 $1,500 $675 ($1.61/100ml) 7折 優惠券 缺貨 免運 ( 52 )'; G=45789.03155092592 }
    @{ Row=5; A=87; C='DORO CAT 豆乳貓 豆腐貓砂, 無味, 7L, 1袋'; D='$115'; E='($1.64/100ml)'; F='DORO CAT 豆乳貓 豆腐貓砂, 無味, 7L, 1袋 特價 46折 $250 $115 ($1.64/100ml) 7折 優惠券 明天 5/13 (二) 預計送達 免運 滿 $490 ( 489 )'; G=45789.03155092592 }
    @{ Row=6; A=86; C='DORO CAT 豆乳貓 添加益生菌與小蘇打顆粒雙重除臭極細豆腐貓砂, 原味, 7L, 5包'; D='$706'; E='($2.02/100ml)'; F='DORO CAT 豆乳貓 添加益生菌與小蘇打顆粒雙重除臭極細豆腐貓砂, 原味, 7L, 5包 特價 57折 $1,250 $706 ($2.02/100ml) 7折 優惠券 明天 5/13 (二) 預計送達 免運 ( 90 )'; G=45789.03155092592 }
    @{ Row=7; A=85; C='DORO CAT 豆乳貓 活性碳與小蘇打雙重除臭顆粒低過敏極細豆腐貓砂, 無香, 7L, 1包'; D='$139'; E='($1.99/100ml)'; F='DORO CAT 豆乳貓 活性碳與小蘇打雙重除臭顆粒低過敏極細豆腐貓砂, 無香, 7L, 1包 特價 56折 $250 $139 ($1.99/100ml) 7折 優惠券 明天 5/13 (二) 預計送達 免運 滿 $490 ( 52 )'; G=45789.03155092592 }
    @{ Row=8; A=84; C='DORO CAT 豆乳貓 添加酵素與小蘇打顆粒雙重消臭極細豆腐貓砂, 木質香, 7L, 4包'; D='$597'; E='($2.13/100ml)'; F='DORO CAT 豆乳貓 添加酵素與小蘇打顆粒雙重消臭極細豆腐貓砂, 木質香, 7L, 4包 特價 6折 $1,000 $597 ($2.13/100ml) 7折 優惠券 明天 5/13 (二) 預計送達 免運 ( 50 )'; G=45789.03155092592 }
    @{ Row=9; A=83; C='DOG CAT STAR 汪喵星球 益生菌消臭條型豆腐砂, 2.7kg, 6袋'; D='$857'; E='($5.29/100g)'; F='DOG CAT STAR 汪喵星球 益生菌消臭條型豆腐砂, 2.7kg, 6袋 特價 41折 $2,100 $857 ($5.29/100g) 7折 優惠券 明天 5/13 (二) 預計送達 免運 ( 2,633 )'; G=45789.03155092592 }
    @{ Row=10; A=82; C='DORO CAT 豆乳貓 混合豆腐礦物貓砂, 無味, 7L, 4袋'; D='$610'; E='($2.18/100ml)'; F='DORO CAT 豆乳貓 混合豆腐礦物貓砂, 無味, 7L, 4袋 特價 61折 $1,000 $610 ($2.18/100ml) 7折 優惠券 明天 5/13 (二) 預計送達 免運 ( 65 )'; G=45789.03155092592 }
    @{ Row=11; A=81; C='DORO CAT 豆乳貓 豆腐貓砂, 抹茶, 7L, 1袋'; D='$114'; E='($1.63/100ml)'; F='DORO CAT 豆乳貓 豆腐貓砂, 抹茶, 7L, 1袋 特價 46折 $250 $114 ($1.63/100ml) 7折 優惠券 明天 5/13 (二) 預計送達 免運 滿 $490 ( 489 )'; G=45789.03155092592 }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in $rows) {
    $row = $r.Row

    # A: id (plain integer)
    $ws.Cells.Item($row, 1).Value = $r.A

    # C: title (text)
    $ws.Cells.Item($row, 3).Value = $r.C

    # D: price - looks like a currency amount ("$298"); force text so Excel
    # doesn't silently reinterpret it as a number, then drop back to the
    # default "Normal" style so no stray number format sticks to the cell.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).Style = "Normal"

    # E: unit_price (text, parenthesised - Excel leaves this alone)
    $ws.Cells.Item($row, 5).Value = $r.E

    # F: full_text (text)
    $ws.Cells.Item($row, 6).Value = $r.F

    # G: timestamp (serial date/time, already styled with the date format)
    $ws.Cells.Item($row, 7).Value = $r.G
}
